$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "PRACTICA Y PRACTICA, MUCHA PRACTICA. 10/09/2024"
$ws.Range("A14").Value = "PRACTICA Y PRACTICA, MUCHA MAS. 10/09/2024"

$ws.Range("A15").Select()
